$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure target cells keep their original text formatting (values such as
# "408.89" or "0.0000111" would otherwise be auto-converted to numbers).
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.618.71'
$ws.Range('E2').Value = '  +6.27%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.294.67'
$ws.Range('E3').Value = '  +1.72%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '408.89'
$ws.Range('E5').Value = '  +3.57%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '111.62'
$ws.Range('E6').Value = '  +4.16%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.287.95'
$ws.Range('E7').Value = '  +1.61%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.563'
$ws.Range('E8').Value = '  -2.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('E9').Value = '  -0.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.613'
$ws.Range('E10').Value = '  -0.83%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.107'
$ws.Range('E11').Value = '  +11.70%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '38.25'
$ws.Range('E12').Value = '  -1.92%  '
$ws.Range('E13').Value = '  -0.03%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.757.55'
$ws.Range('E14').Value = '  +0.15%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '8.09'
$ws.Range('E15').Value = '  -1.13%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.89'
$ws.Range('E16').Value = '  -1.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.269.24'
$ws.Range('E17').Value = '  +0.44%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '60.285.92'
$ws.Range('E18').Value = '  +5.97%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.983'
$ws.Range('E19').Value = '  -4.62%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.39'
$ws.Range('E20').Value = '  -4.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0000111'
$ws.Range('E21').Value = '  +4.89%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.18'
$ws.Range('E22').Value = '  -4.62%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.37'
$ws.Range('E23').Value = '  -4.54%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '293.31'
$ws.Range('E24').Value = '  -0.99%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '72.88'
$ws.Range('E25').Value = '  -1.48%  '
$ws.Range('E26').Value = '  -3.85%  '
$ws.Range('E27').Value = '  +2.71%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '28.42'
$ws.Range('E28').Value = '  +2.14%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.34'
$ws.Range('E29').Value = '  +0.89%  '
$ws.Range('E30').Value = '  +0.55%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.41'
$ws.Range('E31').Value = '  -3.61%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.999'
$ws.Range('E32').Value = '  -0.08%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.08'
$ws.Range('E33').Value = '  -2.85%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.108'
$ws.Range('E34').Value = '  -1.21%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '39.45'
$ws.Range('E35').Value = '  +4.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.38'
$ws.Range('E36').Value = '  +12.62%  '
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0473'
$ws.Range('E37').Value = '  -2.06%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '52.31'
$ws.Range('E38').Value = '  +0.99%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.997'
$ws.Range('E39').Value = '  -0.23%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.06'
$ws.Range('E40').Value = '  +4.15%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.28'
$ws.Range('E41').Value = '  -6.89%  '
$ws.Range('E42').Value = '  +0.03%  '
$ws.Range('E43').Value = '  -1.76%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.86'
$ws.Range('E44').Value = '  -1.34%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.277'
$ws.Range('E45').Value = '  -1.48%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '16.05'
$ws.Range('E46').Value = '  -5.50%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.71'
$ws.Range('E47').Value = '  -6.02%  '
$ws.Range('E48').Value = '  +3.44%  '
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.699.86'
$ws.Range('E49').Value = '  +3.98%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '20.62'
$ws.Range('E50').Value = '  -7.03%  '
$ws.Range('B51').Value = 'Maker'
$ws.Range('C51').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.104.86'
$ws.Range('E51').Value = '  -2.50%  '
